# Auto-generated Excel COM-interop edit script
# Applies numeric value updates (market price recalculation refresh) across all profession sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$edits = @(
  @{ Cell = "H6"; Value = 455 },
  @{ Cell = "I6"; Value = 516 },
  @{ Cell = "K6"; Value = 1548 },
  @{ Cell = "M6"; Value = -1436 },
  @{ Cell = "H12"; Value = 486.2857 },
  @{ Cell = "I12"; Value = 601 },
  @{ Cell = "K12"; Value = 601 },
  @{ Cell = "M12"; Value = -431 },
  @{ Cell = "H17"; Value = 2645.9473 },
  @{ Cell = "J17"; Value = 2645.9473 },
  @{ Cell = "L17"; Value = 7937.841899999999 },
  @{ Cell = "N17"; Value = -8273.841899999999 },
  @{ Cell = "H32"; Value = 4838.5625 },
  @{ Cell = "I32"; Value = 4420.4287 },
  @{ Cell = "J32"; Value = 5163.778 },
  @{ Cell = "K32"; Value = 4420.4287 },
  @{ Cell = "L32"; Value = 5163.778 },
  @{ Cell = "M32"; Value = -4094.4287 },
  @{ Cell = "N32"; Value = -5815.778 },
  @{ Cell = "H33"; Value = 2367.3125 },
  @{ Cell = "I33"; Value = 487.2 },
  @{ Cell = "J33"; Value = 5500.8335 },
  @{ Cell = "K33"; Value = 487.2 },
  @{ Cell = "L33"; Value = 5500.8335 },
  @{ Cell = "M33"; Value = -258.2 },
  @{ Cell = "N33"; Value = -5958.8335 },
  @{ Cell = "H69"; Value = 13776.223 },
  @{ Cell = "I69"; Value = 12855.143 },
  @{ Cell = "J69"; Value = 14362.363 },
  @{ Cell = "K69"; Value = 38565.429 },
  @{ Cell = "L69"; Value = 43087.089 },
  @{ Cell = "M69"; Value = -37691.429 },
  @{ Cell = "N69"; Value = -44835.089 },
  @{ Cell = "H72"; Value = 13776.223 },
  @{ Cell = "I72"; Value = 12855.143 },
  @{ Cell = "J72"; Value = 14362.363 },
  @{ Cell = "K72"; Value = 115696.287 },
  @{ Cell = "L72"; Value = 129261.267 },
  @{ Cell = "M72"; Value = -111328.287 },
  @{ Cell = "N72"; Value = -137997.267 },
  @{ Cell = "H100"; Value = 971.7143 },
  @{ Cell = "I100"; Value = 651.7 },
  @{ Cell = "J100"; Value = 1771.75 },
  @{ Cell = "K100"; Value = 651.7 },
  @{ Cell = "L100"; Value = 1771.75 },
  @{ Cell = "M100"; Value = -110.7 },
  @{ Cell = "N100"; Value = -2853.75 },
  @{ Cell = "H103"; Value = 790.4286 },
  @{ Cell = "I103"; Value = 1899.5 },
  @{ Cell = "J103"; Value = 346.8 },
  @{ Cell = "K103"; Value = 5698.5 },
  @{ Cell = "L103"; Value = 1040.4 },
  @{ Cell = "M103"; Value = -5112.5 },
  @{ Cell = "N103"; Value = -2212.4 },
  @{ Cell = "H107"; Value = 1100 },
  @{ Cell = "I107"; Value = 1022.7778 },
  @{ Cell = "K107"; Value = 1022.7778 },
  @{ Cell = "M107"; Value = 897.2222 },
  @{ Cell = "H112"; Value = 5584 },
  @{ Cell = "I112"; Value = 2329.375 },
  @{ Cell = "J112"; Value = 12093.25 },
  @{ Cell = "K112"; Value = 6988.125 },
  @{ Cell = "L112"; Value = 36279.75 },
  @{ Cell = "M112"; Value = -5880.125 },
  @{ Cell = "N112"; Value = -38495.75 },
  @{ Cell = "H113"; Value = 2113.1 },
  @{ Cell = "I113"; Value = 1904.0869 },
  @{ Cell = "K113"; Value = 1904.0869 },
  @{ Cell = "M113"; Value = 1349.9131 },
  @{ Cell = "H116"; Value = 4281.909 },
  @{ Cell = "I116"; Value = 3991 },
  @{ Cell = "K116"; Value = 3991 },
  @{ Cell = "M116"; Value = -549 },
  @{ Cell = "H129"; Value = 633.7 },
  @{ Cell = "I129"; Value = 633.7 },
  @{ Cell = "K129"; Value = 1901.1 },
  @{ Cell = "M129"; Value = 3098.9 },
  @{ Cell = "H131"; Value = 1810.7778 },
  @{ Cell = "I131"; Value = 1810.7778 },
  @{ Cell = "K131"; Value = 5432.3334 },
  @{ Cell = "M131"; Value = -392.3334000000004 },
  @{ Cell = "H132"; Value = 2016.7 },
  @{ Cell = "I132"; Value = 2016.7 },
  @{ Cell = "K132"; Value = 6050.1 },
  @{ Cell = "M132"; Value = -3520.1 },
  @{ Cell = "H137"; Value = 1633.8889 },
  @{ Cell = "I137"; Value = 1474.4 },
  @{ Cell = "J137"; Value = 2431.3333 },
  @{ Cell = "K137"; Value = 4423.200000000001 },
  @{ Cell = "L137"; Value = 7293.999899999999 },
  @{ Cell = "M137"; Value = -1873.200000000001 },
  @{ Cell = "N137"; Value = -12393.9999 },
  @{ Cell = "H138"; Value = 2945.718 },
  @{ Cell = "J138"; Value = 3079.7693 },
  @{ Cell = "L138"; Value = 9239.3079 },
  @{ Cell = "N138"; Value = -19519.3079 },
  @{ Cell = "H141"; Value = 6643.091 },
  @{ Cell = "I141"; Value = 3845.6667 },
  @{ Cell = "K141"; Value = 11537.0001 },
  @{ Cell = "M141"; Value = -6357.000100000001 }
)
foreach ($e in $edits) { $ws.Range($e.Cell).Value = $e.Value }

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$edits = @(
  @{ Cell = "H28"; Value = 5449.1665 },
  @{ Cell = "I28"; Value = 5449.1665 },
  @{ Cell = "J28"; Value = 0 },
  @{ Cell = "K28"; Value = 5449.1665 },
  @{ Cell = "L28"; Value = 0 },
  @{ Cell = "N28"; Value = -5257.1665 },
  @{ Cell = "H61"; Value = 4608.6 },
  @{ Cell = "I61"; Value = 4332.778 },
  @{ Cell = "K61"; Value = 4332.778 },
  @{ Cell = "M61"; Value = -4120.778 },
  @{ Cell = "H99"; Value = 5449.1665 },
  @{ Cell = "I99"; Value = 5449.1665 },
  @{ Cell = "J99"; Value = 0 },
  @{ Cell = "K99"; Value = 5449.1665 },
  @{ Cell = "L99"; Value = 0 },
  @{ Cell = "N99"; Value = -2454.1665 },
  @{ Cell = "H110"; Value = 1243.875 },
  @{ Cell = "I110"; Value = 1294.6 },
  @{ Cell = "K110"; Value = 1294.6 },
  @{ Cell = "M110"; Value = 750.4000000000001 },
  @{ Cell = "H111"; Value = 60000 },
  @{ Cell = "J111"; Value = 60000 },
  @{ Cell = "L111"; Value = 60000 },
  @{ Cell = "H122"; Value = 1972.3414 },
  @{ Cell = "I122"; Value = 1629.6666 },
  @{ Cell = "K122"; Value = 4888.9998 },
  @{ Cell = "M122"; Value = -2438.9998 },
  @{ Cell = "H132"; Value = 7799.125 },
  @{ Cell = "I132"; Value = 6653.4546 },
  @{ Cell = "K132"; Value = 19960.3638 },
  @{ Cell = "M132"; Value = -17430.3638 },
  @{ Cell = "H136"; Value = 4608.6 },
  @{ Cell = "I136"; Value = 4332.778 },
  @{ Cell = "K136"; Value = 12998.334 },
  @{ Cell = "M136"; Value = -10448.334 },
  @{ Cell = "N111"; Value = -68180 }
)
foreach ($e in $edits) { $ws.Range($e.Cell).Value = $e.Value }
$clears = @("N28", "N99")
foreach ($c in $clears) { $ws.Range($c).ClearContents() }

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$edits = @(
  @{ Cell = "H2"; Value = 48831.668 },
  @{ Cell = "J2"; Value = 48831.668 },
  @{ Cell = "L2"; Value = 48831.668 },
  @{ Cell = "N2"; Value = -49057.668 },
  @{ Cell = "H20"; Value = 25800.2 },
  @{ Cell = "I20"; Value = 3002.6667 },
  @{ Cell = "K20"; Value = 3002.6667 },
  @{ Cell = "M20"; Value = -2755.6667 },
  @{ Cell = "H33"; Value = 45420.668 },
  @{ Cell = "I33"; Value = 32500 },
  @{ Cell = "J33"; Value = 51881 },
  @{ Cell = "K33"; Value = 32500 },
  @{ Cell = "L33"; Value = 51881 },
  @{ Cell = "M33"; Value = -32164 },
  @{ Cell = "N33"; Value = -52553 },
  @{ Cell = "H86"; Value = 3520.84 },
  @{ Cell = "I86"; Value = 1818.6842 },
  @{ Cell = "J86"; Value = 8911 },
  @{ Cell = "K86"; Value = 1818.6842 },
  @{ Cell = "L86"; Value = 8911 },
  @{ Cell = "M86"; Value = -695.6841999999999 },
  @{ Cell = "N86"; Value = -11157 },
  @{ Cell = "H89"; Value = 3520.84 },
  @{ Cell = "I89"; Value = 1818.6842 },
  @{ Cell = "J89"; Value = 8911 },
  @{ Cell = "K89"; Value = 9093.421 },
  @{ Cell = "L89"; Value = 44555 },
  @{ Cell = "M89"; Value = -3477.421 },
  @{ Cell = "N89"; Value = -55787 },
  @{ Cell = "H94"; Value = 969.5454999999999 },
  @{ Cell = "J94"; Value = 1115 },
  @{ Cell = "L94"; Value = 1115 },
  @{ Cell = "N94"; Value = -2017 },
  @{ Cell = "H99"; Value = 1311.2858 },
  @{ Cell = "I99"; Value = 996.5 },
  @{ Cell = "K99"; Value = 996.5 },
  @{ Cell = "M99"; Value = 501.5 }
)
foreach ($e in $edits) { $ws.Range($e.Cell).Value = $e.Value }

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$edits = @(
  @{ Cell = "H7"; Value = 386.83334 },
  @{ Cell = "I7"; Value = 244.94118 },
  @{ Cell = "J7"; Value = 572.38464 },
  @{ Cell = "K7"; Value = 244.94118 },
  @{ Cell = "L7"; Value = 572.38464 },
  @{ Cell = "M7"; Value = -131.94118 },
  @{ Cell = "N7"; Value = -798.38464 },
  @{ Cell = "H22"; Value = 3719.2307 },
  @{ Cell = "I22"; Value = 2233 },
  @{ Cell = "J22"; Value = 4993.143 },
  @{ Cell = "K22"; Value = 2233 },
  @{ Cell = "L22"; Value = 4993.143 },
  @{ Cell = "M22"; Value = -1883 },
  @{ Cell = "N22"; Value = -5693.143 },
  @{ Cell = "H31"; Value = 7407.9736 },
  @{ Cell = "I31"; Value = 3043.4 },
  @{ Cell = "J31"; Value = 15801.385 },
  @{ Cell = "K31"; Value = 3043.4 },
  @{ Cell = "L31"; Value = 15801.385 },
  @{ Cell = "M31"; Value = -2748.4 },
  @{ Cell = "N31"; Value = -16391.385 },
  @{ Cell = "H34"; Value = 7407.9736 },
  @{ Cell = "I34"; Value = 3043.4 },
  @{ Cell = "J34"; Value = 15801.385 },
  @{ Cell = "K34"; Value = 3043.4 },
  @{ Cell = "L34"; Value = 15801.385 },
  @{ Cell = "M34"; Value = -2841.4 },
  @{ Cell = "N34"; Value = -16205.385 },
  @{ Cell = "H50"; Value = 83329.664 },
  @{ Cell = "J50"; Value = 83329.664 },
  @{ Cell = "L50"; Value = 83329.664 },
  @{ Cell = "N50"; Value = -84579.664 },
  @{ Cell = "H51"; Value = 39766.332 },
  @{ Cell = "J51"; Value = 39766.332 },
  @{ Cell = "L51"; Value = 39766.332 },
  @{ Cell = "N51"; Value = -41238.332 },
  @{ Cell = "H57"; Value = 65997.336 },
  @{ Cell = "J57"; Value = 65997.336 },
  @{ Cell = "L57"; Value = 65997.336 },
  @{ Cell = "N57"; Value = -67117.336 },
  @{ Cell = "H59"; Value = 94663 },
  @{ Cell = "J59"; Value = 94663 },
  @{ Cell = "L59"; Value = 94663 },
  @{ Cell = "N59"; Value = -96953 },
  @{ Cell = "H61"; Value = 39766.332 },
  @{ Cell = "J61"; Value = 39766.332 },
  @{ Cell = "L61"; Value = 39766.332 },
  @{ Cell = "N61"; Value = -40462.332 },
  @{ Cell = "H92"; Value = 92631.5 },
  @{ Cell = "J92"; Value = 92631.5 },
  @{ Cell = "L92"; Value = 92631.5 },
  @{ Cell = "N92"; Value = -97623.5 },
  @{ Cell = "H93"; Value = 7933.3335 },
  @{ Cell = "I93"; Value = 7933.3335 },
  @{ Cell = "K93"; Value = 7933.3335 },
  @{ Cell = "H99"; Value = 2208.5 },
  @{ Cell = "I99"; Value = 2263.0833 },
  @{ Cell = "J99"; Value = 2044.75 },
  @{ Cell = "K99"; Value = 2263.0833 },
  @{ Cell = "L99"; Value = 2044.75 },
  @{ Cell = "M99"; Value = -765.0832999999998 },
  @{ Cell = "N99"; Value = -5040.75 },
  @{ Cell = "H107"; Value = 2393.95 },
  @{ Cell = "I107"; Value = 1658 },
  @{ Cell = "J107"; Value = 4601.8 },
  @{ Cell = "K107"; Value = 1658 },
  @{ Cell = "L107"; Value = 4601.8 },
  @{ Cell = "M107"; Value = 262 },
  @{ Cell = "N107"; Value = -8441.799999999999 },
  @{ Cell = "H122"; Value = 3596.5454 },
  @{ Cell = "I122"; Value = 2798.7144 },
  @{ Cell = "K122"; Value = 8396.143199999999 },
  @{ Cell = "M122"; Value = -5946.143199999999 },
  @{ Cell = "H126"; Value = 2208.5 },
  @{ Cell = "I126"; Value = 2263.0833 },
  @{ Cell = "J126"; Value = 2044.75 },
  @{ Cell = "K126"; Value = 6789.249899999999 },
  @{ Cell = "L126"; Value = 6134.25 },
  @{ Cell = "M126"; Value = -4319.249899999999 },
  @{ Cell = "N126"; Value = -11074.25 },
  @{ Cell = "H141"; Value = 298366.5 },
  @{ Cell = "J141"; Value = 379499.78 },
  @{ Cell = "L141"; Value = 379499.78 },
  @{ Cell = "N141"; Value = -389859.78 },
  @{ Cell = "M93"; Value = -6061.3335 }
)
foreach ($e in $edits) { $ws.Range($e.Cell).Value = $e.Value }

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$edits = @(
  @{ Cell = "H2"; Value = 1314.6666 },
  @{ Cell = "J2"; Value = 1299.6666 },
  @{ Cell = "L2"; Value = 7797.9996 },
  @{ Cell = "N2"; Value = -8023.9996 },
  @{ Cell = "H4"; Value = 20139032 },
  @{ Cell = "I4"; Value = 20139032 },
  @{ Cell = "J4"; Value = 0 },
  @{ Cell = "K4"; Value = 60417096 },
  @{ Cell = "L4"; Value = 0 },
  @{ Cell = "N4"; Value = -60416984 },
  @{ Cell = "H8"; Value = 595.4 },
  @{ Cell = "I8"; Value = 595.4 },
  @{ Cell = "K8"; Value = 1786.2 },
  @{ Cell = "M8"; Value = -1647.2 },
  @{ Cell = "H23"; Value = 243.66667 },
  @{ Cell = "I23"; Value = 150 },
  @{ Cell = "J23"; Value = 337.33334 },
  @{ Cell = "K23"; Value = 450 },
  @{ Cell = "L23"; Value = 1012.00002 },
  @{ Cell = "N23"; Value = -1482.00002 },
  @{ Cell = "H24"; Value = 1553.6666 },
  @{ Cell = "J24"; Value = 2231 },
  @{ Cell = "L24"; Value = 6693 },
  @{ Cell = "N24"; Value = -7153 },
  @{ Cell = "H25"; Value = 526 },
  @{ Cell = "I25"; Value = 526 },
  @{ Cell = "K25"; Value = 1578 },
  @{ Cell = "M25"; Value = -1409 },
  @{ Cell = "H30"; Value = 526 },
  @{ Cell = "I30"; Value = 526 },
  @{ Cell = "K30"; Value = 1578 },
  @{ Cell = "M30"; Value = -1476 },
  @{ Cell = "H33"; Value = 898.44446 },
  @{ Cell = "I33"; Value = 1321.5 },
  @{ Cell = "J33"; Value = 52.333332 },
  @{ Cell = "K33"; Value = 7929 },
  @{ Cell = "L33"; Value = 313.999992 },
  @{ Cell = "M33"; Value = -7646 },
  @{ Cell = "N33"; Value = -879.999992 },
  @{ Cell = "H34"; Value = 9586.583000000001 },
  @{ Cell = "I34"; Value = 429.66666 },
  @{ Cell = "J34"; Value = 12638.889 },
  @{ Cell = "K34"; Value = 1288.99998 },
  @{ Cell = "L34"; Value = 37916.667 },
  @{ Cell = "M34"; Value = -1204.99998 },
  @{ Cell = "N34"; Value = -38084.667 },
  @{ Cell = "H39"; Value = 7806 },
  @{ Cell = "I39"; Value = 2500 },
  @{ Cell = "J39"; Value = 9397.799999999999 },
  @{ Cell = "K39"; Value = 7500 },
  @{ Cell = "L39"; Value = 28193.4 },
  @{ Cell = "M39"; Value = -7206 },
  @{ Cell = "N39"; Value = -28781.4 },
  @{ Cell = "H55"; Value = 12186.75 },
  @{ Cell = "J55"; Value = 12642 },
  @{ Cell = "L55"; Value = 37926 },
  @{ Cell = "N55"; Value = -38280 },
  @{ Cell = "H75"; Value = 2812 },
  @{ Cell = "I75"; Value = 475 },
  @{ Cell = "J75"; Value = 3171.5386 },
  @{ Cell = "K75"; Value = 1425 },
  @{ Cell = "L75"; Value = 9514.6158 },
  @{ Cell = "M75"; Value = -427 },
  @{ Cell = "N75"; Value = -11510.6158 },
  @{ Cell = "H78"; Value = 2812 },
  @{ Cell = "I78"; Value = 475 },
  @{ Cell = "J78"; Value = 3171.5386 },
  @{ Cell = "K78"; Value = 4275 },
  @{ Cell = "L78"; Value = 28543.8474 },
  @{ Cell = "M78"; Value = 717 },
  @{ Cell = "N78"; Value = -38527.8474 },
  @{ Cell = "H106"; Value = 9495 },
  @{ Cell = "J106"; Value = 9495 },
  @{ Cell = "L106"; Value = 28485 },
  @{ Cell = "N106"; Value = -30377 },
  @{ Cell = "H113"; Value = 920.3333 },
  @{ Cell = "J113"; Value = 979.9 },
  @{ Cell = "L113"; Value = 2939.7 },
  @{ Cell = "N113"; Value = -7279.7 },
  @{ Cell = "H116"; Value = 8199.75 },
  @{ Cell = "J116"; Value = 7999.5 },
  @{ Cell = "L116"; Value = 23998.5 },
  @{ Cell = "N116"; Value = -30882.5 },
  @{ Cell = "H119"; Value = 5904.5 },
  @{ Cell = "I119"; Value = 3891 },
  @{ Cell = "K119"; Value = 11673 },
  @{ Cell = "M119"; Value = -6835 },
  @{ Cell = "H122"; Value = 899 },
  @{ Cell = "I122"; Value = 299 },
  @{ Cell = "K122"; Value = 2691 },
  @{ Cell = "M122"; Value = -241 },
  @{ Cell = "H129"; Value = 9261733 },
  @{ Cell = "J129"; Value = 13891627 },
  @{ Cell = "L129"; Value = 41674881 },
  @{ Cell = "N129"; Value = -41684881 },
  @{ Cell = "H131"; Value = 3410.3845 },
  @{ Cell = "I131"; Value = 2805 },
  @{ Cell = "K131"; Value = 8415 },
  @{ Cell = "M131"; Value = -3375 },
  @{ Cell = "H139"; Value = 1958.55 },
  @{ Cell = "I139"; Value = 1903.7368 },
  @{ Cell = "K139"; Value = 5711.2104 },
  @{ Cell = "M139"; Value = -571.2103999999999 },
  @{ Cell = "H141"; Value = 7376.1055 },
  @{ Cell = "I141"; Value = 5979.5557 },
  @{ Cell = "K141"; Value = 17938.6671 },
  @{ Cell = "M141"; Value = -12758.6671 },
  @{ Cell = "M23"; Value = -215 }
)
foreach ($e in $edits) { $ws.Range($e.Cell).Value = $e.Value }
$clears = @("N4")
foreach ($c in $clears) { $ws.Range($c).ClearContents() }

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$edits = @(
  @{ Cell = "H70"; Value = 10633 },
  @{ Cell = "I70"; Value = 9900 },
  @{ Cell = "K70"; Value = 9900 },
  @{ Cell = "M70"; Value = -9630 },
  @{ Cell = "H73"; Value = 10633 },
  @{ Cell = "I73"; Value = 9900 },
  @{ Cell = "K73"; Value = 9900 },
  @{ Cell = "M73"; Value = -8964 },
  @{ Cell = "H122"; Value = 4320.8 },
  @{ Cell = "I122"; Value = 1551.1111 },
  @{ Cell = "K122"; Value = 4653.3333 },
  @{ Cell = "M122"; Value = -2203.3333 },
  @{ Cell = "H129"; Value = 70000 },
  @{ Cell = "J129"; Value = 70000 },
  @{ Cell = "L129"; Value = 70000 },
  @{ Cell = "N129"; Value = -80000 },
  @{ Cell = "H132"; Value = 6991.8823 },
  @{ Cell = "I132"; Value = 4187.5835 },
  @{ Cell = "K132"; Value = 12562.7505 },
  @{ Cell = "M132"; Value = -10032.7505 },
  @{ Cell = "H134"; Value = 60730 },
  @{ Cell = "J134"; Value = 60730 },
  @{ Cell = "L134"; Value = 182190 },
  @{ Cell = "N134"; Value = -187260 }
)
foreach ($e in $edits) { $ws.Range($e.Cell).Value = $e.Value }

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$edits = @(
  @{ Cell = "H46"; Value = 1924.72 },
  @{ Cell = "I46"; Value = 881.0909 },
  @{ Cell = "J46"; Value = 2744.7144 },
  @{ Cell = "K46"; Value = 881.0909 },
  @{ Cell = "L46"; Value = 2744.7144 },
  @{ Cell = "M46"; Value = -693.0909 },
  @{ Cell = "N46"; Value = -3120.7144 },
  @{ Cell = "H55"; Value = 790 },
  @{ Cell = "I55"; Value = 201.72223 },
  @{ Cell = "K55"; Value = 201.72223 },
  @{ Cell = "M55"; Value = -28.72223 },
  @{ Cell = "H61"; Value = 4313.56 },
  @{ Cell = "I61"; Value = 1692.6875 },
  @{ Cell = "K61"; Value = 1692.6875 },
  @{ Cell = "M61"; Value = -1490.6875 },
  @{ Cell = "H68"; Value = 6058.05 },
  @{ Cell = "I68"; Value = 3432.6155 },
  @{ Cell = "K68"; Value = 3432.6155 },
  @{ Cell = "M68"; Value = -2683.6155 },
  @{ Cell = "H71"; Value = 6058.05 },
  @{ Cell = "I71"; Value = 3432.6155 },
  @{ Cell = "K71"; Value = 17163.0775 },
  @{ Cell = "M71"; Value = -13419.0775 },
  @{ Cell = "H82"; Value = 1349.9584 },
  @{ Cell = "I82"; Value = 926.61536 },
  @{ Cell = "J82"; Value = 1850.2727 },
  @{ Cell = "K82"; Value = 926.61536 },
  @{ Cell = "L82"; Value = 1850.2727 },
  @{ Cell = "M82"; Value = -565.61536 },
  @{ Cell = "N82"; Value = -2572.2727 },
  @{ Cell = "H85"; Value = 1349.9584 },
  @{ Cell = "I85"; Value = 926.61536 },
  @{ Cell = "J85"; Value = 1850.2727 },
  @{ Cell = "K85"; Value = 926.61536 },
  @{ Cell = "L85"; Value = 1850.2727 },
  @{ Cell = "M85"; Value = 321.38464 },
  @{ Cell = "N85"; Value = -4346.2727 },
  @{ Cell = "H93"; Value = 2931.0667 },
  @{ Cell = "I93"; Value = 1853.2222 },
  @{ Cell = "J93"; Value = 4547.8335 },
  @{ Cell = "K93"; Value = 1853.2222 },
  @{ Cell = "L93"; Value = 4547.8335 },
  @{ Cell = "M93"; Value = -605.2221999999999 },
  @{ Cell = "N93"; Value = -7043.8335 },
  @{ Cell = "H100"; Value = 3798.5945 },
  @{ Cell = "I100"; Value = 3333.077 },
  @{ Cell = "K100"; Value = 3333.077 },
  @{ Cell = "M100"; Value = -2792.077 },
  @{ Cell = "H113"; Value = 4313.56 },
  @{ Cell = "I113"; Value = 1692.6875 },
  @{ Cell = "K113"; Value = 1692.6875 },
  @{ Cell = "M113"; Value = 477.3125 },
  @{ Cell = "H122"; Value = 5073.913 },
  @{ Cell = "J122"; Value = 6288.6 },
  @{ Cell = "L122"; Value = 18865.8 },
  @{ Cell = "N122"; Value = -23765.8 }
)
foreach ($e in $edits) { $ws.Range($e.Cell).Value = $e.Value }

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$edits = @(
  @{ Cell = "H96"; Value = 13316.167 },
  @{ Cell = "I96"; Value = 4900 },
  @{ Cell = "K96"; Value = 4900 },
  @{ Cell = "M96"; Value = -3527 },
  @{ Cell = "H126"; Value = 3063.4 },
  @{ Cell = "I126"; Value = 3152.2144 },
  @{ Cell = "J126"; Value = 1820 },
  @{ Cell = "K126"; Value = 9456.643199999999 },
  @{ Cell = "L126"; Value = 5460 },
  @{ Cell = "M126"; Value = -6986.643199999999 },
  @{ Cell = "N126"; Value = -10400 },
  @{ Cell = "H132"; Value = 3550.6875 },
  @{ Cell = "I132"; Value = 2587.4 },
  @{ Cell = "K132"; Value = 7762.200000000001 },
  @{ Cell = "M132"; Value = -5232.200000000001 },
  @{ Cell = "H136"; Value = 7247.5557 },
  @{ Cell = "I136"; Value = 5104.2144 },
  @{ Cell = "K136"; Value = 15312.6432 },
  @{ Cell = "M136"; Value = -12762.6432 }
)
foreach ($e in $edits) { $ws.Range($e.Cell).Value = $e.Value }

